# Append the new daily-push row (row 72) to Sheet1, matching the existing
# "date / weekday / hour / ranking" table that already fills rows 2-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a literal "yyyy/mm/dd" text label (not a real Excel date),
# just like every other row in the sheet. Force text so Excel's date
# auto-detection doesn't coerce "2025/10/07" into a date serial, then clear
# the formatting override again so the cell keeps the sheet's default style
# (same as all its neighbours, which also carry no explicit style).
$ws.Range("A72").NumberFormat = "@"
$ws.Range("A72").Value = "2025/10/07"
$ws.Range("A72").ClearFormats()

$ws.Range("B72").Value = "火"
$ws.Range("C72").Value = 6
$ws.Range("D72").Value = 6
